$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Split the run containing `"Tsk tsk, you don't trust me?..."` into three
#    runs - "Tsk ", "tsk", ", you don't trust...".  Word itself does this
#    kind of run split when it places spell-check `proofErr` markers around
#    a flagged word; we force the same split point here by briefly anchoring
#    a bookmark on the `tsk` word (bookmarks force a run boundary at their
#    start/end) and then deleting the bookmark again, which leaves the run
#    split in place without leaving the bookmark behind.
# ---------------------------------------------------------------------------
$tskWord = $d.Range(1336, 1339)
$d.Bookmarks.Add("_tmpTskSplit", $tskWord)
$d.Bookmarks.Item("_tmpTskSplit").Delete()

# ---------------------------------------------------------------------------
# 2) Font change: every paragraph that still used Calibri switches to
#    Bookerly (ascii / hAnsi / cs).  The final paragraph already uses
#    Bookerly, so it is left untouched.
# ---------------------------------------------------------------------------
for ($i = 1; $i -le 6; $i++) {
    $pr = $d.Paragraphs.Item($i).Range
    $pr.Font.NameAscii = "Bookerly"
    $pr.Font.NameOther = "Bookerly"
    $pr.Font.NameBi = "Bookerly"
}

# ---------------------------------------------------------------------------
# 3) Move the `_GoBack` bookmark so that it spans the whole document instead
#    of sitting as a zero-width mark in the middle of the last paragraph.
# ---------------------------------------------------------------------------
$d.Bookmarks.Item("_GoBack").Delete()
$wholeDoc = $d.Range(0, $d.Content.End)
$d.Bookmarks.Add("_GoBack", $wholeDoc)

Write-Host "done"
